# Automatic update of files.
# The four observation rows 20-23 are reordered/re-identified:
#   - Row 20 keeps its "Knärot" record but is re-labelled with the Id/coords
#     that used to belong to row 21.
#   - Row 21 becomes the "Blodticka" record that used to be row 23 (including
#     its extra substrate columns J/K/N/AF/AJ/AK/AO).
#   - Row 22 keeps its "Knärot" record but is re-labelled with the Id/coords
#     that used to belong to row 20.
#   - Row 23 becomes the "Knärot" record that used to be row 21 (and loses
#     the substrate columns that belonged to the old "Blodticka" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 20: only Id (A) and coordinates (Q,R) change ----
$ws.Range("A20").Value = 111661832
$ws.Range("Q20").Value = 432076.4609239195
$ws.Range("R20").Value = 6419682.500295377

# ---- Row 21: becomes the "Blodticka" (Meruliopsis taxicola) record ----
$ws.Range("A21").Value = 111661840
$ws.Range("B21").Value = 89793
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 4217
$ws.Range("F21").Value = "Blodticka"
$ws.Range("G21").Value = "Meruliopsis taxicola"
$ws.Range("H21").Value = "(Pers.:Fr.) Bondartsev"
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("Q21").Value = 431768.994999051
$ws.Range("R21").Value = 6419728.1081824
$ws.Range("AF21").Value = ""
$ws.Range("AJ21").Value = "tall"
$ws.Range("AK21").Value = "Pinus sylvestris"
$ws.Range("AO21").Value = "Pinus sylvestris"

# ---- Row 22: only Id (A) and coordinates (Q,R) change ----
$ws.Range("A22").Value = 111661831
$ws.Range("Q22").Value = 432080.3854477856
$ws.Range("R22").Value = 6419662.773410858

# ---- Row 23: becomes the "Knärot" (Goodyera repens) record ----
$ws.Range("A23").Value = 111661838
$ws.Range("B23").Value = 96348
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 220787
$ws.Range("F23").Value = "Knärot"
$ws.Range("G23").Value = "Goodyera repens"
$ws.Range("H23").Value = "(L.) R. Br."
$ws.Range("J23").ClearContents()
$ws.Range("K23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("Q23").Value = 431799.2483237319
$ws.Range("R23").Value = 6419691.460736625
$ws.Range("AF23").ClearContents()
$ws.Range("AJ23").ClearContents()
$ws.Range("AK23").ClearContents()
$ws.Range("AO23").ClearContents()
